$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp string
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 18:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1436123
$ws.Range("C4").Value = 5775
$ws.Range("D4").Value = 310833
$ws.Range("E4").Value = 1039842
$ws.Range("G4").Value = 251
$ws.Range("H4").Value = 85448

# Row 9 - Brasil
$ws.Range("B9").Value = 196158
$ws.Range("C9").Value = 7001
$ws.Range("E9").Value = 104183
$ws.Range("G9").Value = 393
$ws.Range("H9").Value = 13551

# Row 11 - Alemania
$ws.Range("B11").Value = 174284
$ws.Range("C11").Value = 186
$ws.Range("E11").Value = 16116
$ws.Range("G11").Value = 7
$ws.Range("H11").Value = 7868

# Row 22 - Chile
$ws.Range("F22").Value = 663

# Row 25 - Suiza
$ws.Range("E25").Value = 1491
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 1872

# Row 34 - Polonia
$ws.Range("B34").Value = 17615
$ws.Range("C34").Value = 411
$ws.Range("E34").Value = 10036
$ws.Range("G34").Value = 22
$ws.Range("H34").Value = 883

# Row 51 - Chequia
$ws.Range("B51").Value = 8330
$ws.Range("C51").Value = 61
$ws.Range("D51").Value = 5227
$ws.Range("E51").Value = 2811
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 292

# Row 57 - Argelia
$ws.Range("B57").Value = 6442
$ws.Range("C57").Value = 189
$ws.Range("D57").Value = 3158
$ws.Range("E57").Value = 2755
$ws.Range("G57").Value = 7
$ws.Range("H57").Value = 529

# Row 59 - Barein
$ws.Range("B59").Value = 6069
$ws.Range("C59").Value = 253
$ws.Range("D59").Value = 2220
$ws.Range("E59").Value = 3839

# Row 70 - Irak
$ws.Range("B70").Value = 3143
$ws.Range("C70").Value = 111
$ws.Range("D70").Value = 2028
$ws.Range("E70").Value = 1000

# Row 106 - Republica de Chipre
$ws.Range("B106").Value = 907
$ws.Range("C106").Value = 2
$ws.Range("E106").Value = 441

# Row 122 - Jordania
$ws.Range("B122").Value = 586
$ws.Range("C122").Value = 4
$ws.Range("D122").Value = 393
$ws.Range("E122").Value = 184

$wb.Save()
